# Add a new row (row 23) to the Question List worksheet for
# "Binary Tree Upside Down" (leetcode 156), following the same
# pattern as the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the new row's cells
$ws.Range("A23").Value = 22
$ws.Range("B23").Value = "Binary Tree Upside Down"
$ws.Range("D23").Value = "Tree"
$ws.Range("E23").Value = "medium"
$ws.Range("F23").Value = "leetcode 156"

# Match the alignment/style used by the other rows:
#   column A/D/E/F -> centered, column B -> left aligned
$ws.Range("A23").HorizontalAlignment = -4108
$ws.Range("B23").HorizontalAlignment = -4131
$ws.Range("D23").HorizontalAlignment = -4108
$ws.Range("E23").HorizontalAlignment = -4108
$ws.Range("F23").HorizontalAlignment = -4108

# Move the active selection the way Excel would after entering this
# row of data (down to the next empty row in column B).
$null = $ws.Range("B29").Select()
